$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Mc / Industria de las hamburguesas / Mundial / 5 / Comida Rapida / Rudgar Leiva
$ws.Range("A2").Value = "Mc"
$ws.Range("B2").Value = "Industria de las hamburguesas"
$ws.Range("C2").Value = "Mundial"
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = "Comida Rapida"
$ws.Range("F2").Value = "Rudgar Leiva"

# Row 3: BK / Hamburguesas / Mundial / 2 / Comida Rapida / Pedro
$ws.Range("A3").Value = "BK"
$ws.Range("B3").Value = "Hamburguesas"
$ws.Range("C3").Value = "Mundial"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "Comida Rapida"
$ws.Range("F3").Value = "Pedro"

# Row 4: Zara / Venta de ropa / Mundial / 2 / Moda / Juan
$ws.Range("A4").Value = "Zara"
$ws.Range("B4").Value = "Venta de ropa"
$ws.Range("C4").Value = "Mundial"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "Moda"
$ws.Range("F4").Value = "Juan"

# Row 5: Zara / Venta de ropa / Mundial / 2 / Moda / Juan
$ws.Range("A5").Value = "Zara"
$ws.Range("B5").Value = "Venta de ropa"
$ws.Range("C5").Value = "Mundial"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "Moda"
$ws.Range("F5").Value = "Juan"

# Remove the old row 6 (Taco bell) entirely, shrinking the table to 5 rows
$ws.Rows(6).Delete()
